# Reorganize test cases in "Gestion de Activos" sheet and correct CP_GESACT_004.
# - Consolidates CP_GESACT_001..007 into CP_GESACT_001..004
# - Removes now-redundant rows 6-8
# - Updates descriptions / steps / expected & actual results
# - Adjusts row heights and the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = @'
Acceso al la vista Gestión de Activos
'@

$ws.Range("B3").Value = @'
Filtrar ont por estado en "FAILED"
'@
$ws.Range("F3").Value = @'
1.Clic en botón "Seleccionar entidad"
2.Clic en botón "Siguiente"
3.Seleccionar "ont"
4.Clic en botón "Siguiente"
5.Seleccionar fila con ID 9 "FAILED"
6.Hacer clic en el botón "FINALIZAR"
'@
$ws.Range("H3").Value = @'
El sistema registra la selección y finaliza el proceso mostrando la barra de progreso y completando la acción.
'@
$ws.Range("I3").Value = @'
Se seleccionó el registro con estado FAILED y la finalización se ejecutó correctamente.
'@

$ws.Range("B4").Value = @'
actualizacion de estado de la ont FAILED a LOST
'@
$ws.Range("F4").Value = @'
11. Seleccionar primer registro de la tabla (capturar FACTORYSERIAL).
12. Clic en Actualizar estado operativo.
13. Abrir lista de estados.
14. Seleccionar LOST.
15. Diligenciar comentario “test automatización”.
16. Clic en Guardar.
'@
$ws.Range("G4").Value = @'
Serial ONT válido
'@
$ws.Range("H4").Value = @'
El estado operativo del dispositivo se actualiza a LOST y se cierra el modal.
'@
$ws.Range("I4").Value = @'
El estado operativo fue actualizado a LOST y el modal se cerró correctamente.
'@

$ws.Range("B5").Value = @'
Validar actualización de ont a estado LOST
'@
$ws.Range("E5").Value = @'
Haber actualizado el estado en el caso anterior.
'@
$ws.Range("F5").Value = @'
17–23. Repetir selección de entidad “elemento secundario”, tipo “ONT”, seleccionar fila con ID 10 “LOST”, clic en FINALIZAR y esperar la barra de progreso.
'@
$ws.Range("G5").Value = @'
Registro con estado LOST
'@
$ws.Range("H5").Value = @'
El sistema completa nuevamente el flujo para el dispositivo con estado LOST.
'@
$ws.Range("I5").Value = @'
El flujo se completó correctamente para el registro con estado LOST.
'@

# Remove the now-merged/obsolete rows (former CP_GESACT_005/006/007)
$ws.Rows("6:8").Delete()

# Row heights for the remaining detail rows
$ws.Rows(3).RowHeight = 148.5
$ws.Rows(4).RowHeight = 171
$ws.Rows(5).RowHeight = 171

# Restore the printed page orientation
$ws.PageSetup.Orientation = 1

# Move the active selection to A5 (top-left scrolled to row 4)
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("A5").Select()
